$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for columns D, M, N, O, P, S
$d2 = $ws.Range("D2").Value2
$m2 = $ws.Range("M2").Value2
$n2 = $ws.Range("N2").Value2
$o2 = $ws.Range("O2").Value2
$p2 = $ws.Range("P2").Value2
$s2 = $ws.Range("S2").Value2

$d3 = $ws.Range("D3").Value2
$m3 = $ws.Range("M3").Value2
$n3 = $ws.Range("N3").Value2
$o3 = $ws.Range("O3").Value2
$p3 = $ws.Range("P3").Value2
$s3 = $ws.Range("S3").Value2

$ws.Range("D2").Value2 = $d3
$ws.Range("M2").Value2 = $m3
$ws.Range("N2").Value2 = $n3
$ws.Range("O2").Value2 = $o3
$ws.Range("P2").Value2 = $p3
$ws.Range("S2").Value2 = $s3

$ws.Range("D3").Value2 = $d2
$ws.Range("M3").Value2 = $m2
$ws.Range("N3").Value2 = $n2
$ws.Range("O3").Value2 = $o2
$ws.Range("P3").Value2 = $p2
$ws.Range("S3").Value2 = $s2
